$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# Paragraph 9 (1-based): "Ideally: we will have ~26 raters rating ~26 papers each (~10-15 hours)"
#   -> "Ideally: we will have ~20 raters rating ~30 papers each"
$tr.Paragraphs(9,1).Text = "Ideally: we will have ~20 raters rating ~30 papers each"

# Paragraph 11 (1-based): "+ meetings, discussion, collaborative analysis & consensus-building…"
#   -> "Develop the TMS-RAT in stages, reach consensus, input to analysis & writing"
$tr.Paragraphs(11,1).Text = "Develop the TMS-RAT in stages, reach consensus, input to analysis & writing"

# Paragraph 13 (1-based): "*** Co-authorship of TMS-RAT 🐀 paper ***" (first/last run underlined)
#   -> "(and co-authorship of TMS-RAT 🐀 paper!)" (no underline)
$para13 = $tr.Paragraphs(13,1)
$run1 = $para13.Runs(1,1)
$run1.Text = "(and co-authorship of TMS-RAT "
$run1.Font.Underline = $false

$para13 = $tr.Paragraphs(13,1)
$run3 = $para13.Runs(3,1)
$run3.Text = "paper!)"
$run3.Font.Underline = $false
